# Generate Report for Handback
# Refreshes the "Latest HO Xliff Generate Date" / handoff & handback
# datetime stamps that get regenerated each time the handback report runs.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# ffdaea81... row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-24 23:05:38"

# zh-cn handback sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) for the ffdaea81... row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-24 23:05:33"
$wsZhCn.Range("K2").Value = "2016-08-24 23:05:50"

# de-de handback sheet: same two columns for the ffdaea81... row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-24 23:05:38"
$wsDeDe.Range("K2").Value = "2016-08-24 23:05:57"
